$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("movimentos")

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "café"
$ws.Cells.Item(3, 3).Value = "ENTRADA"
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = "2026-01-20 11:31:29"
